$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the narrative text for the remaining patient rows
$ws.Range("D2").Value = "Patient presents with 3-month history of depressed mood"
$ws.Range("E2").Value = "Patient is a 35-year-old male with major depressive disorder"

$ws.Range("D3").Value = "Patient reports increasing anxiety over past 6 months"
$ws.Range("E3").Value = "Patient is a 28-year-old female with generalized anxiety disorder"

$ws.Range("D4").Value = "Patient describes alternating periods of high and low mood"
$ws.Range("E4").Value = "Patient is a 42-year-old male with bipolar disorder"

# Drop the last two test patients (Sarah Williams, Robert Brown) entirely
$ws.Range("A5:F6").EntireRow.Delete()
